# Fix bug: update vehicle categories for Budget
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 3 new rows right after row 23 (current last row of the
#    "Budget" block) to make room for 3 additional vehicle categories.
#    This pushes everything from the old row 24 onward down by 3 rows.
# ------------------------------------------------------------------
$ws.Rows("24:26").Insert()

# ------------------------------------------------------------------
# 2. Populate the 3 new rows with the new Budget vehicle categories,
#    and fix the existing "Premium UTE" label (row 23) -> "Premium
#    Ute". The values below are written in this particular order so
#    that the new shared-string table entries line up with the
#    original author's edit order.
# ------------------------------------------------------------------
$ws.Range("B24").Value2 = 9
$ws.Range("B25").Value2 = ""
$ws.Range("B26").Value2 = 8

# Row 25: Compact Hybrid (also has a new_id / new_category_name pair)
$ws.Range("C25").Value2 = "Compact Hybrid"
$ws.Range("D25").Value2 = "Toyota Corolla Hybrid or similar"
$ws.Range("G25").Value2 = "Compact Hybrid Car"
$ws.Range("F25").Value2 = 23

# Row 24: Premium Minivan
$ws.Range("C24").Value2 = "Premium Minivan"
$ws.Range("D24").Value2 = "Hyundai iMax or similar"

# Fix row 23's category label
$ws.Range("C23").Value2 = "Premium Ute"

# Row 26: Standard SUV
$ws.Range("C26").Value2 = "Standard SUV"
$ws.Range("D26").Value2 = "Toyota Rav4 2WD or similar"

# ------------------------------------------------------------------
# 4. Formatting for the new cells, matching the surrounding rows.
#    - Column A / B already inherited the row-23 look (center/center)
#      from the insert; merge the whole "Budget" block (A13:A26) so it
#      matches the rest of the category column.
#    - F25 should use the same green fill used by the other "new_id"
#      helper cells in column F (copy from F22).
# ------------------------------------------------------------------
$ws.Range("A13:A26").Merge()

$ws.Range("F22").Copy()
$ws.Range("F25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A24").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A24").VerticalAlignment = -4108
$ws.Range("A25").HorizontalAlignment = -4108
$ws.Range("A25").VerticalAlignment = -4108
$ws.Range("A26").HorizontalAlignment = -4108
$ws.Range("A26").VerticalAlignment = -4108

$ws.Range("B24").HorizontalAlignment = -4108
$ws.Range("B24").VerticalAlignment = -4108
$ws.Range("B25").HorizontalAlignment = -4108
$ws.Range("B25").VerticalAlignment = -4108
$ws.Range("B26").HorizontalAlignment = -4108
$ws.Range("B26").VerticalAlignment = -4108

# ------------------------------------------------------------------
# 5. Update the selected cell to reflect where the author ended up.
# ------------------------------------------------------------------
$ws.Range("F25").Select()
